$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style()
    $cell.NumberFormat = "@"
    $cell.Value = "'" + $val
    $cell.Style = $origStyle
}

Set-TextCell $ws "D2" '43.669.50'
Set-TextCell $ws "E2" '  -0.21%  '

Set-TextCell $ws "D3" '2.324.92'
Set-TextCell $ws "E3" '  +4.33%  '

Set-TextCell $ws "E4" '  +0.13%  '

Set-TextCell $ws "D5" '270.92'
Set-TextCell $ws "E5" '  -1.20%  '

Set-TextCell $ws "D6" '94.82'
Set-TextCell $ws "E6" '  +8.01%  '

Set-TextCell $ws "D7" '0.626'
Set-TextCell $ws "E7" '  +1.85%  '

Set-TextCell $ws "D9" '0.618'
Set-TextCell $ws "E9" '  +2.30%  '

Set-TextCell $ws "D10" '44.63'
Set-TextCell $ws "E10" '  -0.66%  '

Set-TextCell $ws "D11" '0.0941'
Set-TextCell $ws "E11" '  +2.32%  '

Set-TextCell $ws "D12" '7.98'
Set-TextCell $ws "E12" '  +4.00%  '

Set-TextCell $ws "E13" '  +0.21%  '

Set-TextCell $ws "D14" '2.680.49'
Set-TextCell $ws "E14" '  +4.78%  '

Set-TextCell $ws "D15" '15.70'
Set-TextCell $ws "E15" '  +4.78%  '

Set-TextCell $ws "D16" '0.854'
Set-TextCell $ws "E16" '  +7.95%  '

Set-TextCell $ws "D17" '2.341.34'
Set-TextCell $ws "E17" '  +5.14%  '

Set-TextCell $ws "D18" '43.657.99'
Set-TextCell $ws "E18" '  -0.05%  '

Set-TextCell $ws "E19" '  +3.34%  '

Set-TextCell $ws "D20" '6.33'
Set-TextCell $ws "E20" '  +6.18%  '

Set-TextCell $ws "D21" '71.91'
Set-TextCell $ws "E21" '  +2.38%  '

Set-TextCell $ws "D22" '237.83'
Set-TextCell $ws "E22" '  +2.35%  '

Set-TextCell $ws "D23" '2.24'
Set-TextCell $ws "E23" '  -4.96%  '

Set-TextCell $ws "D24" '9.54'
Set-TextCell $ws "E24" '  +8.99%  '

Set-TextCell $ws "E25" '  -0.11%  '

Set-TextCell $ws "D26" '2.54'
Set-TextCell $ws "E26" '  -1.38%  '

Set-TextCell $ws "D27" '11.29'
Set-TextCell $ws "E27" '  +4.43%  '

Set-TextCell $ws "D28" '3.42'
Set-TextCell $ws "E28" '  -1.88%  '

Set-TextCell $ws "E29" '  -0.81%  '

Set-TextCell $ws "D30" '38.41'
Set-TextCell $ws "E30" '  -1.53%  '

Set-TextCell $ws "E31" '  +8.36%  '

Set-TextCell $ws "D32" '172.15'
Set-TextCell $ws "E32" '  -0.37%  '

Set-TextCell $ws "D33" '0.0892'
Set-TextCell $ws "E33" '  -1.23%  '

Set-TextCell $ws "D34" '5.45'
Set-TextCell $ws "E34" '  +1.62%  '

Set-TextCell $ws "D35" '0.126'
Set-TextCell $ws "E35" '  +2.68%  '

Set-TextCell $ws "D36" '0.0355'
Set-TextCell $ws "E36" '  +0.40%  '

Set-TextCell $ws "E37" '  -3.56%  '

Set-TextCell $ws "D38" '4.34'
Set-TextCell $ws "E38" '  +1.74%  '

Set-TextCell $ws "D39" '3.38'
Set-TextCell $ws "E39" '  -1.18%  '

Set-TextCell $ws "D40" '2.36'
Set-TextCell $ws "E40" '  +8.73%  '

Set-TextCell $ws "D41" '0.232'
Set-TextCell $ws "E41" '  +11.86%  '

Set-TextCell $ws "D42" '1.37'
Set-TextCell $ws "E42" '  +21.16%  '

Set-TextCell $ws "D43" '11.96'
Set-TextCell $ws "E43" '  -4.05%  '

Set-TextCell $ws "B44" 'FraxShare'
Set-TextCell $ws "C44" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell $ws "D44" '9.07'
Set-TextCell $ws "E44" '  +7.14%  '

Set-TextCell $ws "B45" 'MultiversX'
Set-TextCell $ws "C45" 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextCell $ws "D45" '61.43'
Set-TextCell $ws "E45" '  -3.75%  '

Set-TextCell $ws "D46" '5.35'
Set-TextCell $ws "E46" '  -0.64%  '

Set-TextCell $ws "E47" '  +4.25%  '

Set-TextCell $ws "D48" '100.46'
Set-TextCell $ws "E48" '  +0.08%  '

Set-TextCell $ws "D49" '1.22'
Set-TextCell $ws "E49" '  +2.04%  '

Set-TextCell $ws "D50" '2.556.62'
Set-TextCell $ws "E50" '  +4.62%  '

Set-TextCell $ws "D51" '0.181'
Set-TextCell $ws "E51" '  +12.53%  '
